$d = $word.ActiveDocument

# The document ends with (after the "LOQ4053..." requirements paragraph):
#   [empty paragraph]
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#    pages. Original theme under Creative Commons Attribution"
#   [empty paragraph]
#   [page-break paragraph]
#
# The edit removes the first three of those (the blank line plus the two
# "Ver no Jupiter" / copyright lines), so "LOQ4053..." is immediately
# followed by the remaining blank paragraph and then the page break.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $startPara = $d.Paragraphs.Item($target - 1)   # blank paragraph right before it
    $endPara   = $d.Paragraphs.Item($target + 1)    # copyright paragraph right after it

    $r = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $r.Delete()
}
